$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 with "Test Exp 15" experiment details
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Test Exp 15"
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 0.45
$ws.Range("E16").Value = "Local"
$ws.Range("F16").Value = -1
$ws.Range("G16").Value = "28*28"
$ws.Range("H16").Value = "32*32"
$ws.Range("I16").Value = "3,4,5"

# Match the formatting of the other data rows (style index 2: left aligned)
$ws.Range("A16:H16").HorizontalAlignment = -4131

# Update the selected cell as recorded in the workbook view
$ws.Range("E19").Select()
